$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: row 7 becomes the (now empty) single summary line, qty 0 ---
$ws.Range("A7").Value = 0
$ws.Range("C7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("Q7").Value = ""

# --- Step 2: carry the "totals" row (old row 12, P12:Q12) up onto row 8 ---
$ws.Range("P12:Q12").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8:Q8").Merge()

# --- Step 3: drop the old row-8 product entry (now superseded by the totals cells) ---
$ws.Range("A8:B8").UnMerge()
$ws.Range("C8:G8").UnMerge()
$ws.Range("H8:K8").UnMerge()
$ws.Range("L8:M8").UnMerge()
$ws.Range("N8:O8").UnMerge()
$ws.Range("A8:O8").Clear()

# --- Step 4: refresh the footer timestamp (still on its original row, row 13, at this point) ---
$ws.Range("A13").Value = "Wednesday, 17 September, 2025 12:45 PM"

# --- Step 5: remove the now-superfluous product rows (old rows 9-11) and the old totals row (old row 12) ---
$ws.Rows("9:12").Delete()
